# Chiffres COVID-19 Valais - daily update
# - One extra positive case recorded on row 259 (2020-10-14): C259 425 -> 426
#   (cascades through the "Cumul cas positifs" running-total formula in column B)
# - Row 367 (2021-02-11) new-case count corrected 68 -> 67 so the running
#   total (column B) returns to its original trajectory from that date on.
# - Rows 433-437 (2021-05-03 .. 2021-05-07): corrected new-case counts, and
#   row 437 also gets updated SI/intubation figures.
# - Rows 438-440 (2021-05-08 .. 2021-05-10): new daily figures entered
#   (previously blank placeholder rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One more positive case on 2020-10-14 (row 259)
$ws.Range("C259").Value = 426

# Correction on 2021-02-11 (row 367)
$ws.Range("C367").Value = 67

# Corrections 2021-05-03 .. 2021-05-06 (rows 433-436)
$ws.Range("C433").Value = 110
$ws.Range("C434").Value = 86
$ws.Range("C435").Value = 61
$ws.Range("C436").Value = 60

# 2021-05-07 (row 437): revised case count plus updated SI figures
$ws.Range("C437").Value = 76
$ws.Range("E437").Value = 8
$ws.Range("F437").Value = 6

# 2021-05-08 (row 438): newly entered daily figures
$ws.Range("C438").Value = 37
$ws.Range("E438").Value = 8
$ws.Range("F438").Value = 6
$ws.Range("G438").Value = 18

# 2021-05-09 (row 439): newly entered daily figures
$ws.Range("C439").Value = 26
$ws.Range("E439").Value = 8
$ws.Range("F439").Value = 6
$ws.Range("G439").Value = 19

# 2021-05-10 (row 440): newly entered daily figures
$ws.Range("C440").Value = 7
$ws.Range("E440").Value = 8
$ws.Range("F440").Value = 6
$ws.Range("G440").Value = 21

# Columns L (intubés) and M (sorties) on these rows are formatted as Text
# (numFmtId 49, "@"); a direct .Value = 0 would be stored as the text "0"
# instead of a number. Flip the format to General for the write, then
# restore it to Text so the style index (s="17" / s="18") is unchanged but
# the stored value is a real number, matching how the rest of the sheet's
# data (e.g. L433:M437) is stored.
foreach ($addr in @("L438", "M438", "L439", "M439", "L440", "M440")) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Value = 0
    $ws.Range($addr).NumberFormat = "@"
}
